# Add all subclasses to content list
# - Subclasses sheet: 8 new rows (71-78) with new subclasses + source-doc hyperlinks
# - Feats sheet: 1 new row (17) for the "Warcaller" content-dependant feat + hyperlink

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Subclasses sheet
# ---------------------------------------------------------------------------
$sub = $wb.Worksheets.Item("Subclasses")

# row, Name, Base Class, Revised Subclass, Source Doc (cell text), hyperlink URL,
# hyperlink display text (only set when it differs from the cell text),
# Development Status, Release Status, Added To Subclass Sheet?, Supporting Content Status
$subRows = @(
    @(71, "College of Revelry",  "Bard",      "No",  "Bard College - College of Revely",   "https://editor.gmbinder.com/documents/edit/-N8RevelryCollege0", $null,                              "Playtest Ready",  "Not Released", "Yes", "None"),
    @(72, "College of Pacts",    "Bard",      "No",  "Bard College - College of Pacts",     "https://editor.gmbinder.com/documents/edit/-N8PactsCollege000", $null,                              "Playtest Ready",  "Not Released", "Yes", "None"),
    @(73, "College of Choir",    "Bard",      "No",  "Bard College - College of Choir",     "https://editor.gmbinder.com/documents/edit/-N8ChoirCollege000", "Bard College - Collge of Choir",   "Needs Clean Up",  "Not Released", "No",  "None"),
    @(74, "College of Finality", "Bard",      "No",  "Bard College - College of Finality",  "https://editor.gmbinder.com/documents/edit/-N8FinalityCollege0", $null,                              "Needs Clean Up",  "Not Released", "No",  "None"),
    @(75, "Divine Domain - War", "Cleric",    "Yes", "Divine Domain - War",                 "https://editor.gmbinder.com/documents/edit/-N8DivineDomainWar0", $null,                              "Needs Clean Up",  "Not Released", "No",  "None"),
    @(76, "Demon Soul",          "Barbarian", "No",  "Path of the Demon Soul",              "https://editor.gmbinder.com/documents/edit/-N8DemonSoulPath000", "Primal Path - Demon Soul",         "Playtest Ready",  "Not Released", "Yes", "Unknown"),
    @(77, "Warcaller",           "Barbarian", "No",  "Path of the Warcaller",               "https://editor.gmbinder.com/documents/edit/-N8WarcallerPath000", $null,                              "Playtest Ready",  "Not Released", "Yes", "Unknown"),
    @(78, "Storm Herald",        "Barbarian", "Yes", "Path of the Storm Herald",            "https://editor.gmbinder.com/documents/edit/-N8StormHeraldPath0", $null,                              "Playtest Ready",  "Not Released", "Yes", "None")
)

foreach ($r in $subRows) {
    $rowNum = $r[0]

    # Start each new row off looking like the row above it (font/alignment/number format/etc.)
    $sub.Range("A" + ($rowNum - 1) + ":H" + ($rowNum - 1)).Copy() | Out-Null
    $sub.Range("A" + $rowNum + ":H" + $rowNum).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $sub.Range("A$rowNum").Value = $r[1]
    $sub.Range("B$rowNum").Value = $r[2]
    $sub.Range("C$rowNum").Value = $r[3]
    $sub.Range("D$rowNum").Value = $r[4]
    $sub.Range("E$rowNum").Value = $r[7]
    $sub.Range("F$rowNum").Value = $r[8]
    $sub.Range("G$rowNum").Value = $r[9]
    $sub.Range("H$rowNum").Value = $r[10]

    $displayText = $r[6]
    if ($null -eq $displayText) {
        $sub.Hyperlinks.Add($sub.Range("D$rowNum"), $r[5]) | Out-Null
    } else {
        $sub.Hyperlinks.Add($sub.Range("D$rowNum"), $r[5], "", "", $displayText) | Out-Null
    }
    # Hyperlinks.Add can stomp the cell's displayed text with the link's display
    # text; re-assert the real Source Doc text so the cell itself always shows
    # the correct (non-typo'd) string even when the hyperlink's cached display
    # text is stale.
    $sub.Range("D$rowNum").Value = $r[4]

    # Re-apply the D-column format (Hyperlinks.Add nudges the cell style) so it keeps
    # looking like the rest of the Source Doc column.
    $sub.Range("D" + ($rowNum - 1)).Copy() | Out-Null
    $sub.Range("D$rowNum").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

$excel.CutCopyMode = $false

$sub.Activate()
$sub.Range("E80").Select() | Out-Null

# ---------------------------------------------------------------------------
# Feats sheet
# ---------------------------------------------------------------------------
$feats = $wb.Worksheets.Item("Feats")

$feats.Range("A16:G16").Copy() | Out-Null
$feats.Range("A17:G17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$feats.Range("A17").Value = "Warcaller"
$feats.Range("B17").Value = "Content Dependant"
$feats.Range("C17").Value = "None"
$feats.Range("D17").Value = "No"
$feats.Range("E17").Value = "Path of the Warcaller"
$feats.Range("F17").Value = "Playtest Ready"
$feats.Range("G17").Value = "Not Released"

$feats.Hyperlinks.Add($feats.Range("E17"), "https://editor.gmbinder.com/documents/edit/-N8WarcallerPath000") | Out-Null
$feats.Range("E17").Value = "Path of the Warcaller"

$feats.Range("E16").Copy() | Out-Null
$feats.Range("E17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = $false

$feats.Columns.Item(2).ColumnWidth = 18.7109375

$feats.Activate()
$feats.Range("B20").Select() | Out-Null

$sub.Activate()
